$d = $word.ActiveDocument

# 1. Update "Curso (semestre ideal): EB (6)" -> "EB (8)"
$d.Content.Find.Execute("Curso (semestre ideal): EB (6)", $false, $false, $false, $false, $false, $true, 1, $false, "Curso (semestre ideal): EB (8)", 2)

# 2. Insert a new ListBullet paragraph "4873328 - Fernando Segato" right after
#    the "Docente(s) Responsável(eis)" heading paragraph.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Docente(s) Responsável(eis) ") {
        $p.Range.InsertParagraphAfter()
        $newPara = $p.Next()
        $newPara.Range.Text = "4873328 - Fernando Segato"
        $newPara.Style = "ListBullet"
        break
    }
}

# 3 & 4. Update requisitos text. Both strings live as sibling <w:r> runs inside
# the same paragraph; replacing them individually via Find/Range.Text causes
# the engine to coalesce the (identically-formatted) runs into one <w:r> with
# two <w:t>/<w:br/> pairs. Rebuild the whole paragraph body via InsertXML
# (scoped to the paragraph's text, excluding its trailing paragraph mark) so
# the two runs stay distinct, matching the source structure.
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text
    if ($txt -like "*LOT2040*" -and $txt -like "*LOT2053*") {
        $target = $d.Range($p.Range.Start, $p.Range.End - 1)
        $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
          '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
          '<pkg:xmlData>' +
          '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
          '<w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr>' +
          '<w:r><w:t>LOT2040 -  Engenharia Genética Teórica e Prática  (Requisito fraco)</w:t><w:br/></w:r>' +
          '<w:r><w:t>LOT2053 -  Microbiologia: da Teoria à Prática  (Requisito fraco)</w:t><w:br/></w:r>' +
          '</w:p>' +
          '</w:body>' +
          '</w:document>' +
          '</pkg:xmlData></pkg:part></pkg:package>'
        $target.InsertXML($xml)
        break
    }
}
